$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Force the Price column (which often looks numeric, e.g. "610.74")
    # to stay a plain text cell, matching the source inline-string XML,
    # by prefixing with an apostrophe (Excel's "store as text" marker)
    # and then resetting the style so no quotePrefix/number-format is
    # left behind on the cell.
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "66.729.73"
$ws.Range("E2").Value = "  +0.99%  "

# Row 3 - Ethereum
Set-TextCell "D3" "3.634.72"
$ws.Range("E3").Value = "  +2.00%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-TextCell "D5" "610.74"
$ws.Range("E5").Value = "  +0.65%  "

# Row 6 - Solana
Set-TextCell "D6" "150.11"
$ws.Range("E6").Value = "  +3.44%  "

# Row 7 - LidoStakedEther
Set-TextCell "D7" "3.634.27"
$ws.Range("E7").Value = "  +2.02%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.09%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.20%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.40%  "

# Row 11 - Toncoin
Set-TextCell "D11" "7.95"
$ws.Range("E11").Value = "  +0.15%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +1.32%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextCell "D13" "4.251.06"
$ws.Range("E13").Value = "  +2.01%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  +1.26%  "

# Row 15 - Avalanche
Set-TextCell "D15" "30.00"
$ws.Range("E15").Value = "  +0.10%  "

# Row 16 - WrappedEther
Set-TextCell "D16" "3.632.91"
$ws.Range("E16").Value = "  +1.81%  "

# Row 17 - WrappedBTC
Set-TextCell "D17" "66.827.48"
$ws.Range("E17").Value = "  +0.97%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  +1.49%  "

# Row 19 - Uniswap
Set-TextCell "D19" "11.64"
$ws.Range("E19").Value = "  +2.29%  "

# Row 20 - Polkadot
Set-TextCell "D20" "6.38"
$ws.Range("E20").Value = "  +2.84%  "

# Row 21 - Chainlink
Set-TextCell "D21" "15.17"
$ws.Range("E21").Value = "  +2.18%  "

# Row 22 - BitcoinCash
Set-TextCell "D22" "429.08"
$ws.Range("E22").Value = "  -0.05%  "

# Row 23 - Polygon
Set-TextCell "D23" "0.622"
$ws.Range("E23").Value = "  +1.41%  "

# Row 24 - Litecoin
Set-TextCell "D24" "78.94"
$ws.Range("E24").Value = "  -0.25%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  -0.01%  "

# Row 26 - PEPE
Set-TextCell "D26" "0.0000124"
$ws.Range("E26").Value = "  +4.89%  "

# Row 27 - RenderToken
Set-TextCell "D27" "8.41"
$ws.Range("E27").Value = "  +6.00%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("E28").Value = "  +5.21%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  +0.83%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  -0.13%  "

# Row 31 - RenzoRestakedETH
Set-TextCell "D31" "3.632.91"
$ws.Range("E31").Value = "  +2.15%  "

# Row 32 - Fetch.AI
$ws.Range("E32").Value = "  +1.38%  "

# Row 33 - Kaspa
$ws.Range("E33").Value = "  +4.20%  "

# Row 34 - EthereumClassic
Set-TextCell "D34" "25.51"
$ws.Range("E34").Value = "  -0.22%  "

# Row 35 - Aptos
$ws.Range("E35").Value = "  +0.82%  "

# Row 37 - NEARProtocol
Set-TextCell "D37" "5.70"
$ws.Range("E37").Value = "  +1.68%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -1.16%  "

# Row 39 - Monero
Set-TextCell "D39" "176.79"
$ws.Range("E39").Value = "  +1.01%  "

# Row 40 - Hedera
Set-TextCell "D40" "0.0864"
$ws.Range("E40").Value = "  +2.05%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  +1.16%  "

# Row 42 - Mantle
$ws.Range("E42").Value = "  +0.71%  "

# Row 43 - Stacks
$ws.Range("E43").Value = "  -1.19%  "

# Row 44 - OKB
$ws.Range("E44").Value = "  -0.57%  "

# Row 45 - dogwifhat
$ws.Range("E45").Value = "  +8.20%  "

# Row 46 - FirstDigitalUSD
$ws.Range("E46").Value = "  +0.00%  "

# Row 47 - was InjectiveProtocol, now ONDO
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextCell "D47" "1.18"
$ws.Range("E47").Value = "  -1.32%  "

# Row 48 - was ONDO, now InjectiveProtocol
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D48" "25.15"
$ws.Range("E48").Value = "  -2.36%  "

# Row 49 - EnergySwap
Set-TextCell "D49" "24.09"
$ws.Range("E49").Value = "  +2.72%  "

# Row 50 - Cosmos
Set-TextCell "D50" "7.22"
$ws.Range("E50").Value = "  +1.40%  "

# Row 51 - SuiNetwork
$ws.Range("E51").Value = "  +1.97%  "
